$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'67.793.18"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Formula = "'3.809.78"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Formula = "'597.99"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Formula = "'167.84"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Formula = "'0.161"
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("D10").Formula = "'6.30"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Formula = "'0.449"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Formula = "'36.05"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Formula = "'4.447.59"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Formula = "'3.826.39"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("E16").Value = "  +5.37%  "
$ws.Range("D17").Formula = "'67.783.54"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Formula = "'461.91"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Formula = "'9.96"
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Formula = "'0.0000156"
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("D24").Formula = "'83.63"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Formula = "'12.11"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("D26").Formula = "'2.11"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Formula = "'10.03"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Formula = "'3.958.00"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +4.69%  "
$ws.Range("D32").Formula = "'7.29"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").Formula = "'29.71"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Formula = "'9.10"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Formula = "'0.998"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Formula = "'3.749.05"
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").Formula = "'0.1000"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").Formula = "'0.999"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Formula = "'48.20"
$ws.Range("E44").Value = "  +2.86%  "
$ws.Range("D45").Formula = "'43.77"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").Formula = "'149.18"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("D48").Formula = "'8.33"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Formula = "'397.53"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").Formula = "'26.90"
$ws.Range("E51").Value = "  +6.99%  "
